$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 5 and 6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1056
$ws1.Range("F6").Value = 2448

# Sheet "全部类型" (All types) - rows 7 and 8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1056
$ws4.Range("F8").Value = 2448
